$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fill in the "Final Deliverable" scores (column G) for rows 2-5
$ws.Range("G2").Value = 25
$ws.Range("G3").Value = 25
$ws.Range("G4").Value = 25
$ws.Range("G5").Value = 25

# Update the active cell selection to G6 (matching the authored edit)
$ws.Range("G6").Select()
